# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2023-09-10 (serial 45179) to 2023-09-11 (serial 45180).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows span from row 2 through row 205.
$firstRow = 2
$lastRow = 205

$range = $ws.Range("C$firstRow`:C$lastRow")
$range.Value = 45180
